# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (DAMSLTag, DialogAct)
$updates = @{
    12 = @("sd", "Statement-non-opinion")
    24 = @("sd", "Statement-non-opinion")
    25 = @("sd", "Statement-non-opinion")
    26 = @("sd", "Statement-non-opinion")
    27 = @("sd", "Statement-non-opinion")
    39 = @("aa", "Agree/Accept")
    40 = @("aa", "Agree/Accept")
    60 = @("sd", "Statement-non-opinion")
    66 = @("sd", "Statement-non-opinion")
    69 = @("sv", "Statement-opinion")
    71 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
